# Daily attendance processing - 2025-11-07 07:22:35
# Normalizes the "Recorded By" (column G) lists so that the literal
# entry "System" is moved to the end of the comma-separated list
# (preserving the relative order of the remaining entries), and - for
# the rare rows where no "System" entry is present - the recorder
# names are placed in alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row on the sheet so we cover every data row.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ([string]::IsNullOrEmpty($value)) { continue }

    # Split the comma separated recorder list, trimming whitespace.
    $rawParts = $value -split ','
    $parts = @()
    foreach ($rp in $rawParts) {
        $parts += $rp.Trim()
    }

    if ($parts.Length -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals('System')) { $hasSystem = $true }
    }

    if ($hasSystem) {
        # Move the first exact-case "System" entry to the end,
        # keeping every other entry (e.g. a lowercase "system") in
        # its original relative position.
        $newParts = @()
        $removed = $false
        foreach ($p in $parts) {
            if ((-not $removed) -and $p.Equals('System')) {
                $removed = $true
                continue
            }
            $newParts += $p
        }
        $newParts += 'System'
        $newValue = [string]::Join(', ', $newParts)
    }
    else {
        # No "System" entry present - order the recorders alphabetically.
        $sortedParts = $parts | Sort-Object
        $newValue = [string]::Join(', ', $sortedParts)
    }

    if (-not $newValue.Equals($value)) {
        $cell.Value = $newValue
    }
}
